$d = $word.ActiveDocument

# The original tables were marked with <w:bidi w:val="0"/> inside <w:tblPr>.
# Re-assigning TableDirection (even to the same logical "LTR" state) makes
# Word re-emit this as the modern <w:bidiVisual w:val="0"/> element instead
# of the legacy <w:bidi w:val="0"/> element, for every table in the body.
foreach ($t in $d.Tables) {
    $t.TableDirection = 0
}

# Append the extra sentence to the paragraph about the wi-fi receiver.
$old = "Dependendo do local talvez seja preciso que o receptor wi-fi seja mais robusto que possua maior alcance de sinal, sendo isso possível com o uso de uma antena externa, além de fornecer suporte para cartão de memória."
$new = "Dependendo do local talvez seja preciso que o receptor wi-fi seja mais robusto que possua maior alcance de sinal, sendo isso possível com o uso de uma antena externa, além de fornecer suporte para cartão de memória. Porém seu preço é bem elevado."

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

Write-Output "done"
